# RPA datasets push 2023-11-07
# Update row 17 (캡스톤파트너스): confirmed offer price and offering amount

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D17: 확정공모가 (confirmed offering price) changes from "-" to 4000
# Stored as text (matches sibling cells D18:D21, which hold "2000", "25000",
# etc. as text), so force a text number format before assigning the
# numeric-looking string, then restore the cell style to the sheet default
# so the cell keeps its original (unstyled) appearance.
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4000"
$ws.Range("D17").Style = "Normal"

# E17: 공모금액(백만) (offering amount, millions) changes from 5107 to 6384
$ws.Range("E17").Value = 6384
